# Update '想去人数' (attendee interest counts) column F across all sheets
# to match the regenerated gh-pages data snapshot (commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 9615
$ws.Range("F9").Value = 38
$ws.Range("F10").Value = 714
$ws.Range("F11").Value = 2118
$ws.Range("F13").Value = 1631
$ws.Range("F14").Value = 2737
$ws.Range("F15").Value = 138
$ws.Range("F16").Value = 4082
$ws.Range("F17").Value = 338
$ws.Range("F18").Value = 163
$ws.Range("F20").Value = 221
$ws.Range("F21").Value = 242
$ws.Range("F22").Value = 32
$ws.Range("F24").Value = 82
$ws.Range("F25").Value = 280
$ws.Range("F26").Value = 3876
$ws.Range("F27").Value = 5
$ws.Range("F28").Value = 3366
$ws.Range("F29").Value = 1103
$ws.Range("F30").Value = 200
$ws.Range("F31").Value = 492
$ws.Range("F32").Value = 4337
$ws.Range("F33").Value = 74
$ws.Range("F34").Value = 292
$ws.Range("F35").Value = 403
$ws.Range("F36").Value = 271

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 23

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 202
$ws.Range("F3").Value = 999

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 202
$ws.Range("F4").Value = 999
$ws.Range("F9").Value = 9615
$ws.Range("F11").Value = 38
$ws.Range("F12").Value = 714
$ws.Range("F13").Value = 2118
$ws.Range("F15").Value = 1631
$ws.Range("F17").Value = 2737
$ws.Range("F18").Value = 138
$ws.Range("F19").Value = 4082
$ws.Range("F20").Value = 338
$ws.Range("F21").Value = 163
$ws.Range("F23").Value = 221
$ws.Range("F24").Value = 242
$ws.Range("F25").Value = 32
$ws.Range("F26").Value = 23
$ws.Range("F28").Value = 82
$ws.Range("F29").Value = 280
$ws.Range("F30").Value = 3876
$ws.Range("F31").Value = 5
$ws.Range("F32").Value = 3366
$ws.Range("F33").Value = 1103
$ws.Range("F34").Value = 200
$ws.Range("F35").Value = 492
$ws.Range("F36").Value = 4337
$ws.Range("F37").Value = 74
$ws.Range("F38").Value = 292
$ws.Range("F39").Value = 403
$ws.Range("F40").Value = 271

